$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task")
$ws.Range("A1").Value = "test"
